$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 26. This shifts the former rows 26-29
#    down to 27-30 (and copies formatting/styles down, same as Excel does
#    when a row is inserted).
# ---------------------------------------------------------------------------
$ws.Rows("26:26").Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new row 26 with the "Recurso: Bodega de tiempos" entry.
#    Columns A (ID) and B (ID_Padre) are left blank, as in the source diff.
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = "Recurso: Bodega de tiempos"
$ws.Range("E26").Value = "archivo"
$ws.Range("F26").Value = "historico de bodega de tiempos para limpieza"
$ws.Range("G26").Value = "https://grupomarval.sharepoint.com/:f:/r/sites/Lean-BIM/Documentos%20compartidos/LEAN/13.%20Cloud%20Data/0.%20ETL%20Colab/13.0.1.%20Historico%20Bodegas%20de%20tiempo?csf=1&web=1&e=y0SGaf"

# ---------------------------------------------------------------------------
# 3. Fix up the hyperlinks. Inserting the row moved the underlying cells of
#    the hyperlinks that used to sit on G26 and G28 down to G27 and G29, but
#    the hyperlink anchors themselves stay bound to the old addresses, so we
#    recreate them pointing at the correct (shifted) cells, and add a new
#    hyperlink for the newly created G26 cell.
# ---------------------------------------------------------------------------
# Delete matching hyperlinks one at a time, re-scanning the live collection
# after every removal (deleting while holding onto stale references from a
# pre-built list can skip entries once the collection re-indexes).
$changed = $true
while ($changed) {
    $changed = $false
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$G$26' -or $addr -eq '$G$28') {
            $h.Delete()
            $changed = $true
            break
        }
    }
}

$ws.Hyperlinks.Add($ws.Range("G27"), "https://colab.research.google.com/drive/1XhnifKohqfhBAwTACpzS7d6D4qcx_sl0", "scrollTo=3ZOFGD83yd0f")
$ws.Hyperlinks.Add($ws.Range("G29"), "https://console.cloud.google.com/bigquery?project=modelomarval&ws=!1m4!1m3!3m2!1smodelomarval!2sproyectos")
$ws.Hyperlinks.Add($ws.Range("G26"), "https://grupomarval.sharepoint.com/:f:/r/sites/Lean-BIM/Documentos%20compartidos/LEAN/13.%20Cloud%20Data/0.%20ETL%20Colab/13.0.1.%20Historico%20Bodegas%20de%20tiempo?csf=1&web=1&e=y0SGaf")

# Note: the new row 26 already inherited the correct "Hipervinculo"-like
# style (centered fill + underlined hyperlink font) on G26 automatically
# from the row-insert operation (Excel copies the format of the row above
# when inserting), matching the other URL cells in the column, so no
# explicit re-styling is necessary here.

# ---------------------------------------------------------------------------
# 5. Grow the table (ListObject) so it covers the new row.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H30"))

# ---------------------------------------------------------------------------
# 6. Widen column G to fit the long new URL.
# ---------------------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 253.67

# ---------------------------------------------------------------------------
# 7. Update the view: the sheet was scrolled so column G / row 10 is the
#    top-left visible cell, and the active selection moved to G12.
# ---------------------------------------------------------------------------
$win = $excel.Windows.Item(1)
$win.ScrollRow = 10
$win.ScrollColumn = 7
$ws.Range("G12").Select()
